$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy format of the last existing data row (row 21) down to new rows 22-28
# so new rows match column A bold/border/centered style and plain number cells elsewhere.
$ws.Range("A21:AD21").Copy()
$ws.Range("A22:AD28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 22
$ws.Range("A22").Value2 = 20
$ws.Range("B22").Value2 = 239
$ws.Range("C22").Value2 = 20
$ws.Range("D22").Value2 = 0.369921875
$ws.Range("E22").Value2 = -0.03078125
$ws.Range("F22").Value2 = 57.74231311258944
$ws.Range("G22").Value2 = 0.04928864760245783
$ws.Range("H22").Value2 = 0.001907359417408593
$ws.Range("I22").Value2 = -1
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 0.437142534032867
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = 0
$ws.Range("N22").Value2 = 15.1720875
$ws.Range("O22").Value2 = 1.905537500000001
$ws.Range("P22").Value2 = 12.791425
$ws.Range("Q22").Value2 = 2
$ws.Range("R22").Value2 = 80783
$ws.Range("S22").Value2 = 1
$ws.Range("T22").Value2 = $true
$ws.Range("U22").Value2 = 77
$ws.Range("V22").Value2 = 10
$ws.Range("W22").Value2 = 240
$ws.Range("X22").Value2 = 302
$ws.Range("Y22").Value2 = 2
$ws.Range("Z22").Value2 = 82
$ws.Range("AA22").Value2 = "atq"
$ws.Range("AB22").Value2 = 1
$ws.Range("AC22").Value2 = "traditional"
$ws.Range("AD22").Value2 = 1.567242433984148

# Row 23
$ws.Range("A23").Value2 = 21
$ws.Range("B23").Value2 = 239
$ws.Range("C23").Value2 = 21
$ws.Range("D23").Value2 = 0.378671875
$ws.Range("E23").Value2 = -0.029140625
$ws.Range("F23").Value2 = 59.00248231300491
$ws.Range("G23").Value2 = 0.04955877145819353
$ws.Range("H23").Value2 = 0.001911389520358967
$ws.Range("I23").Value2 = -1
$ws.Range("J23").Value2 = 0
$ws.Range("K23").Value2 = 0.4540006788732852
$ws.Range("L23").Value2 = 0
$ws.Range("M23").Value2 = 0
$ws.Range("N23").Value2 = 15.1737125
$ws.Range("O23").Value2 = 1.980875000000001
$ws.Range("P23").Value2 = 12.79205
$ws.Range("Q23").Value2 = 2
$ws.Range("R23").Value2 = 80520
$ws.Range("S23").Value2 = 1
$ws.Range("T23").Value2 = $true
$ws.Range("U23").Value2 = 77
$ws.Range("V23").Value2 = 10
$ws.Range("W23").Value2 = 240
$ws.Range("X23").Value2 = 302
$ws.Range("Y23").Value2 = 2
$ws.Range("Z23").Value2 = 82
$ws.Range("AA23").Value2 = "atq"
$ws.Range("AB23").Value2 = 1
$ws.Range("AC23").Value2 = "traditional"
$ws.Range("AD23").Value2 = 1.629682359345268

# Row 24
$ws.Range("A24").Value2 = 22
$ws.Range("B24").Value2 = 239
$ws.Range("C24").Value2 = 22
$ws.Range("D24").Value2 = 0.365078125
$ws.Range("E24").Value2 = -0.032890625
$ws.Range("F24").Value2 = 60.0986020400894
$ws.Range("G24").Value2 = 0.0497999872003506
$ws.Range("H24").Value2 = 0.001916193157425531
$ws.Range("I24").Value2 = -1
$ws.Range("J24").Value2 = 0
$ws.Range("K24").Value2 = 0.4708585412821902
$ws.Range("L24").Value2 = 0
$ws.Range("M24").Value2 = 0
$ws.Range("N24").Value2 = 15.17575625
$ws.Range("O24").Value2 = 2.056212500000001
$ws.Range("P24").Value2 = 12.79265
$ws.Range("Q24").Value2 = 2
$ws.Range("R24").Value2 = 80388
$ws.Range("S24").Value2 = 1
$ws.Range("T24").Value2 = $true
$ws.Range("U24").Value2 = 77
$ws.Range("V24").Value2 = 10
$ws.Range("W24").Value2 = 240
$ws.Range("X24").Value2 = 302
$ws.Range("Y24").Value2 = 2
$ws.Range("Z24").Value2 = 82
$ws.Range("AA24").Value2 = "atq"
$ws.Range("AB24").Value2 = 1
$ws.Range("AC24").Value2 = "traditional"
$ws.Range("AD24").Value2 = 1.646271365061731

# Row 25
$ws.Range("A25").Value2 = 23
$ws.Range("B25").Value2 = 239
$ws.Range("C25").Value2 = 23
$ws.Range("D25").Value2 = 0.3721875
$ws.Range("E25").Value2 = -0.0309375
$ws.Range("F25").Value2 = 61.43230919474978
$ws.Range("G25").Value2 = 0.05000641473902703
$ws.Range("H25").Value2 = 0.001916285444844785
$ws.Range("I25").Value2 = -1
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 0.4877137686888458
$ws.Range("L25").Value2 = 0
$ws.Range("M25").Value2 = 0
$ws.Range("N25").Value2 = 15.17824375
$ws.Range("O25").Value2 = 2.1315375
$ws.Range("P25").Value2 = 12.7932375
$ws.Range("Q25").Value2 = 2
$ws.Range("R25").Value2 = 80277
$ws.Range("S25").Value2 = 1
$ws.Range("T25").Value2 = $true
$ws.Range("U25").Value2 = 77
$ws.Range("V25").Value2 = 10
$ws.Range("W25").Value2 = 240
$ws.Range("X25").Value2 = 302
$ws.Range("Y25").Value2 = 2
$ws.Range("Z25").Value2 = 82
$ws.Range("AA25").Value2 = "atq"
$ws.Range("AB25").Value2 = 1
$ws.Range("AC25").Value2 = "traditional"
$ws.Range("AD25").Value2 = 1.690213258276134

# Row 26
$ws.Range("A26").Value2 = 24
$ws.Range("B26").Value2 = 239
$ws.Range("C26").Value2 = 24
$ws.Range("D26").Value2 = 0.3696875
$ws.Range("E26").Value2 = -0.03078125
$ws.Range("F26").Value2 = 62.73743421508154
$ws.Range("G26").Value2 = 0.05010663951524311
$ws.Range("H26").Value2 = 0.001915920986484906
$ws.Range("I26").Value2 = -1
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 0.5045701660188114
$ws.Range("L26").Value2 = 0
$ws.Range("M26").Value2 = 0
$ws.Range("N26").Value2 = 15.181175
$ws.Range("O26").Value2 = 2.206875
$ws.Range("P26").Value2 = 12.793825
$ws.Range("Q26").Value2 = 2
$ws.Range("R26").Value2 = 80491
$ws.Range("S26").Value2 = 1
$ws.Range("T26").Value2 = $true
$ws.Range("U26").Value2 = 77
$ws.Range("V26").Value2 = 10
$ws.Range("W26").Value2 = 240
$ws.Range("X26").Value2 = 302
$ws.Range("Y26").Value2 = 2
$ws.Range("Z26").Value2 = 82
$ws.Range("AA26").Value2 = "atq"
$ws.Range("AB26").Value2 = 1
$ws.Range("AC26").Value2 = "traditional"
$ws.Range("AD26").Value2 = 1.696000534247075

# Row 27
$ws.Range("A27").Value2 = 25
$ws.Range("B27").Value2 = 239
$ws.Range("C27").Value2 = 25
$ws.Range("D27").Value2 = 0.37265625
$ws.Range("E27").Value2 = -0.031328125
$ws.Range("F27").Value2 = 63.90097069853066
$ws.Range("G27").Value2 = 0.05033477273181355
$ws.Range("H27").Value2 = 0.001923591010199174
$ws.Range("I27").Value2 = -1
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 0.5214279424411044
$ws.Range("L27").Value2 = 0
$ws.Range("M27").Value2 = 0
$ws.Range("N27").Value2 = 15.1845875
$ws.Range("O27").Value2 = 2.2822125
$ws.Range("P27").Value2 = 12.7945
$ws.Range("Q27").Value2 = 2
$ws.Range("R27").Value2 = 80448
$ws.Range("S27").Value2 = 1
$ws.Range("T27").Value2 = $true
$ws.Range("U27").Value2 = 77
$ws.Range("V27").Value2 = 10
$ws.Range("W27").Value2 = 240
$ws.Range("X27").Value2 = 302
$ws.Range("Y27").Value2 = 2
$ws.Range("Z27").Value2 = 82
$ws.Range("AA27").Value2 = "atq"
$ws.Range("AB27").Value2 = 1
$ws.Range("AC27").Value2 = "traditional"
$ws.Range("AD27").Value2 = 1.762057682353713

# Row 28
$ws.Range("A28").Value2 = 26
$ws.Range("B28").Value2 = 239
$ws.Range("C28").Value2 = 26
$ws.Range("D28").Value2 = 0.37125
$ws.Range("E28").Value2 = -0.03046875
$ws.Range("F28").Value2 = 64.71992978692553
$ws.Range("G28").Value2 = 0.05044477057173122
$ws.Range("H28").Value2 = 0.00192518646410966
$ws.Range("I28").Value2 = -1
$ws.Range("J28").Value2 = 0
$ws.Range("K28").Value2 = 0.538285347725167
$ws.Range("L28").Value2 = 0
$ws.Range("M28").Value2 = 0
$ws.Range("N28").Value2 = 15.1884375
$ws.Range("O28").Value2 = 2.35755
$ws.Range("P28").Value2 = 12.7951625
$ws.Range("Q28").Value2 = 2
$ws.Range("R28").Value2 = 80492
$ws.Range("S28").Value2 = 1
$ws.Range("T28").Value2 = $true
$ws.Range("U28").Value2 = 77
$ws.Range("V28").Value2 = 10
$ws.Range("W28").Value2 = 240
$ws.Range("X28").Value2 = 302
$ws.Range("Y28").Value2 = 2
$ws.Range("Z28").Value2 = 82
$ws.Range("AA28").Value2 = "atq"
$ws.Range("AB28").Value2 = 1
$ws.Range("AC28").Value2 = "traditional"
$ws.Range("AD28").Value2 = 1.751368406321544
